$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Copy formatting (style) of column G into the new column H for every
#    existing row, so the new column visually matches its neighbour.
$ws.Range("G1:G50").Copy() | Out-Null
$ws.Range("H1:H50").PasteSpecial(-4122) | Out-Null

# 2) Header for the new "Is verified" column.
$ws.Range("H1").Value = "Is verified"

# 3) Per-row "Is verified" values.
#    Row 2 -> Yes, Row 3 -> No, Row 4 -> (blank), Row 5 -> yes, Row 6 -> no
$ws.Range("H2").Value = "Yes"
$ws.Range("H3").Value = "No"
$ws.Range("H5").Value = "yes"
$ws.Range("H6").Value = "no"

# 4) Fix the F5 e-mail text (was split across two rich-text runs "t" +
#    "test@example.com", now a single run "ttest@example.com"). The logical
#    text is already "ttest@example.com" (two runs concatenate to that), so
#    a plain .Value assignment is a same-value no-op that keeps the old
#    2-run split. Toggling the first character's Bold on/off forces the
#    engine to re-derive/merge the rich-text runs into one, collapsing them
#    into a single run once their formatting matches again.
$ch = $ws.Range("F5").Characters(1, 1)
$ch.Font.Bold = $true
$ch.Font.Bold = $false

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("F2"), "mailto:test@example.com", "", "", "test@example.com")
$ws.Hyperlinks.Add($ws.Range("F3"), "mailto:test@example.com", "", "", "test1@example.com")
$ws.Hyperlinks.Add($ws.Range("F4"), "mailto:test@example.com", "", "", "test2@example.com")
$ws.Hyperlinks.Add($ws.Range("F5"), "mailto:test@example.com", "", "", "ttest@example.com")
$ws.Hyperlinks.Add($ws.Range("F6"), "mailto:test@example.com", "", "", "test3@example.com")

# 5) Match the column width used by the other data columns (F, G).
$ws.Columns.Item(8).ColumnWidth = 16.83
